$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The document's "last edit position" bookmark (_GoBack) moves from right
#    after "Four" (start of the title) to inside "(40 min of work)" - between
#    "(40 m" and the rest of the word "min". Re-adding a bookmark with the
#    same name moves it (Word keeps only one bookmark per name), and Word
#    automatically splits the run that spans the insertion point so the
#    bookmark sits cleanly between two runs.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(40 m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $target)

# ---------------------------------------------------------------------------
# 2. Set font to Calibri for non-heading text: the document-wide "Normal"
#    paragraph style (which headings are based on but override) switches
#    from Tahoma 12pt to Calibri 11pt.
# ---------------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 11
